$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197734951972961
$ws.Range("B1").Value = 3.403714895248413
$ws.Range("C1").Value = 2.617115259170532
$ws.Range("D1").Value = 2.412417650222778
$ws.Range("E1").Value = 1.988042593002319
